# Update odds values on Sheet1 for the 2026-01-02 Betfair back/lay workbook.
# These are the individual cell corrections captured in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Melbourne Victory vs Perth Glory
$ws.Range("F2").Value = 1.85
$ws.Range("G2").Value = 1.86
$ws.Range("I2").Value = 4.7
$ws.Range("L2").Value = 1.33
$ws.Range("Q2").Value = 1.75
$ws.Range("T2").Value = 1.72
$ws.Range("W2").Value = 2.16
$ws.Range("X2").Value = 18
$ws.Range("AE2").Value = 55
$ws.Range("AO2").Value = 55

# Row 3 - Al Najma Club vs Al-Khaleej Saihat
$ws.Range("T3").Value = 1.54

# Row 4 - Al-Ettifaq vs Al-Akhdoud
$ws.Range("G4").Value = 1.68
$ws.Range("J4").Value = 4
$ws.Range("N4").Value = 2.64
$ws.Range("P4").Value = 1.96
$ws.Range("Q4").Value = 1.54
$ws.Range("T4").Value = 1.73
$ws.Range("U4").Value = 1.86
$ws.Range("W4").Value = 2.46
$ws.Range("X4").Value = 20
$ws.Range("Y4").Value = 28
$ws.Range("AB4").Value = 11.5
$ws.Range("AC4").Value = 13.5
$ws.Range("AD4").Value = 32
$ws.Range("AF4").Value = 14
$ws.Range("AG4").Value = 14
$ws.Range("AH4").Value = 28
$ws.Range("AJ4").Value = 23
$ws.Range("AK4").Value = 26
$ws.Range("AN4").Value = 14

# Row 5 - Kabylie vs MC Alger
$ws.Range("M5").Value = 1.09
$ws.Range("O5").Value = 1.09

# Row 6 - Omonia FC Aradippou vs Digenis Ypsona
$ws.Range("G6").Value = 2.68
$ws.Range("W6").Value = 1.59

# Row 7 - Al Ahli vs Al Nassr
$ws.Range("H7").Value = 1.53
$ws.Range("N7").Value = 5.8
$ws.Range("O7").Value = 1.16
$ws.Range("Z7").Value = 13.5
$ws.Range("AB7").Value = 990
$ws.Range("AC7").Value = 13.5
$ws.Range("AD7").Value = 11
$ws.Range("AE7").Value = 980
$ws.Range("AF7").Value = 60
$ws.Range("AJ7").Value = 150
$ws.Range("AK7").Value = 60
$ws.Range("AL7").Value = 65
$ws.Range("AM7").Value = 75
$ws.Range("AO7").Value = 6.2

# Row 8 - Gil Vicente vs Sporting Lisbon
$ws.Range("F8").Value = 7.6
$ws.Range("U8").Value = 1.79
$ws.Range("Z8").Value = 8.199999999999999
$ws.Range("AE8").Value = 19
$ws.Range("AF8").Value = 65
$ws.Range("AH8").Value = 28
$ws.Range("AI8").Value = 44
$ws.Range("AK8").Value = 150
$ws.Range("AM8").Value = 170
$ws.Range("AO8").Value = 9.800000000000001

# Row 9 - Eibar vs Mirandes
$ws.Range("N9").Value = 2.54
$ws.Range("Q9").Value = 1.84
$ws.Range("AN9").Value = 20

# Row 10 - Toulouse vs Lens
$ws.Range("P10").Value = 1.92
$ws.Range("Q10").Value = 2.02
$ws.Range("AM10").Value = 95

# Row 11 - Cagliari vs AC Milan
$ws.Range("O11").Value = 1.33

# Row 12 - Rayo Vallecano vs Getafe
$ws.Range("I12").Value = 4.5
$ws.Range("AJ12").Value = 30

# Row 13 - Guimaraes vs CD Nacional Funchal
$ws.Range("M13").Value = 1.08
$ws.Range("N13").Value = 3.25
$ws.Range("Q13").Value = 2.1
$ws.Range("R13").Value = 1.29
$ws.Range("S13").Value = 3.95
$ws.Range("T13").Value = 1.95
$ws.Range("U13").Value = 1.91
$ws.Range("Y13").Value = 16
$ws.Range("AA13").Value = 130
$ws.Range("AB13").Value = 7.6
$ws.Range("AD13").Value = 20
$ws.Range("AG13").Value = 10
$ws.Range("AH13").Value = 22
$ws.Range("AL13").Value = 42
$ws.Range("AN13").Value = 16
$ws.Range("AO13").Value = 100

Write-Host "Applied $($ws.UsedRange.Count) cell updates"
